$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 273, shifting the existing rows 273-329 down to 275-331.
$ws.Rows("273:274").Insert()

# Populate the first new row (273)
$ws.Cells.Item(273, 1).Value2 = 10
$ws.Cells.Item(273, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(273, 3).Value2 = "La Araucanía"
$ws.Cells.Item(273, 4).Value2 = 45211
$ws.Cells.Item(273, 5).Value2 = 9
$ws.Cells.Item(273, 6).Value2 = "Fruta"
$ws.Cells.Item(273, 7).Value2 = 100101
$ws.Cells.Item(273, 8).Value2 = "Berries"
$ws.Cells.Item(273, 9).Value2 = 100112025
$ws.Cells.Item(273, 10).Value2 = "Frutilla"
$ws.Cells.Item(273, 11).Value2 = "Sin especificar"
$ws.Cells.Item(273, 12).Value2 = "Primera"
$ws.Cells.Item(273, 13).Value2 = 465
$ws.Cells.Item(273, 14).Value2 = 13000
$ws.Cells.Item(273, 15).Value2 = 14000
$ws.Cells.Item(273, 16).Value2 = 13398
$ws.Cells.Item(273, 17).Value2 = "$/bandeja 7 kilos"
$ws.Cells.Item(273, 18).Value2 = "Provincia de Melipilla"
$ws.Cells.Item(273, 19).Value2 = 1914
$ws.Cells.Item(273, 20).Value2 = 7

# Populate the second new row (274)
$ws.Cells.Item(274, 1).Value2 = 10
$ws.Cells.Item(274, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(274, 3).Value2 = "La Araucanía"
$ws.Cells.Item(274, 4).Value2 = 45211
$ws.Cells.Item(274, 5).Value2 = 9
$ws.Cells.Item(274, 6).Value2 = "Fruta"
$ws.Cells.Item(274, 7).Value2 = 100101
$ws.Cells.Item(274, 8).Value2 = "Berries"
$ws.Cells.Item(274, 9).Value2 = 100112025
$ws.Cells.Item(274, 10).Value2 = "Frutilla"
$ws.Cells.Item(274, 11).Value2 = "Sin especificar"
$ws.Cells.Item(274, 12).Value2 = "Segunda"
$ws.Cells.Item(274, 13).Value2 = 500
$ws.Cells.Item(274, 14).Value2 = 8000
$ws.Cells.Item(274, 15).Value2 = 9000
$ws.Cells.Item(274, 16).Value2 = 8400
$ws.Cells.Item(274, 17).Value2 = "$/caja 7 kilos"
$ws.Cells.Item(274, 18).Value2 = "Región de La Araucanía"
$ws.Cells.Item(274, 19).Value2 = 1200
$ws.Cells.Item(274, 20).Value2 = 7
